$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before the old "grupo" column (C), shifting
# grupo/rol/NO CAMBIAR TITULOS from C,D,E to F,G,H.
$ws.Columns("C:E").Insert()

# Rename the "punto" header to "punto de encuentro".
$ws.Range("B1").Value = "punto de encuentro"

# New header row cells.
$ws.Range("C1").Value = "presente"
$ws.Range("D1").Value = "pago"
$ws.Range("E1").Value = "estado"

# New helper/comment row cells.
$ws.Range("C2").Value = "[Si / No]"
$ws.Range("D2").Value = "[Si / No]"
$ws.Range("E2").Value = "[Sin Contactar, Sin Interés, Sin Confirmar, Confirmado]"

# New columns get the same width as column B (punto de encuentro).
$ws.Range("C1:E1").ColumnWidth = 64

# Selection moved to B5 in the saved file.
$ws.Range("B5").Select() | Out-Null
